# Apply latest crypto price/volume updates to columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the cell's original style
# (needed for numeric-looking strings like "0.9973" which Excel would
# otherwise silently convert to a Double, changing the stored cell type).
function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "28.981.50"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.828.85"
$ws.Range("E3").Value = "  -0.05%  "

Set-TextValue $ws.Range("D4") "0.9973"
$ws.Range("E4").Value = "  -0.16%  "

Set-TextValue $ws.Range("D5") "243.60"
$ws.Range("E5").Value = "  +1.05%  "

Set-TextValue $ws.Range("D6") "0.6320"
$ws.Range("E6").Value = "  +1.41%  "

Set-TextValue $ws.Range("D7") "0.9981"
$ws.Range("E7").Value = "  -0.22%  "

Set-TextValue $ws.Range("D8") "0.07516"
$ws.Range("E8").Value = "  -0.54%  "

Set-TextValue $ws.Range("D9") "0.2938"
$ws.Range("E9").Value = "  +0.75%  "

Set-TextValue $ws.Range("D10") "22.92"
$ws.Range("E10").Value = "  +0.67%  "

Set-TextValue $ws.Range("D11") "0.07723"
$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").Value = "1.835.47"
$ws.Range("E12").Value = "  +0.25%  "

Set-TextValue $ws.Range("D13") "4.993"
$ws.Range("E13").Value = "  +0.71%  "

Set-TextValue $ws.Range("D14") "0.6706"
$ws.Range("E14").Value = "  +0.90%  "

Set-TextValue $ws.Range("D15") "83.10"
$ws.Range("E15").Value = "  +1.03%  "

Set-TextValue $ws.Range("D16") "0.000009722"
$ws.Range("E16").Value = "  +7.62%  "

Set-TextValue $ws.Range("D17") "6.074"
$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").Value = "29.020.52"
$ws.Range("E18").Value = "  -0.06%  "

Set-TextValue $ws.Range("D19") "12.54"
$ws.Range("E19").Value = "  +1.67%  "

Set-TextValue $ws.Range("D20") "226.65"
$ws.Range("E20").Value = "  +0.83%  "

Set-TextValue $ws.Range("D21") "0.9976"

Set-TextValue $ws.Range("D22") "7.178"
$ws.Range("E22").Value = "  -0.34%  "

Set-TextValue $ws.Range("D23") "0.9985"
$ws.Range("E23").Value = "  -0.22%  "

Set-TextValue $ws.Range("D24") "159.65"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  +3.30%  "

Set-TextValue $ws.Range("D26") "8.540"
$ws.Range("E26").Value = "  +1.80%  "

Set-TextValue $ws.Range("D27") "17.90"
$ws.Range("E27").Value = "  +0.37%  "

Set-TextValue $ws.Range("D28") "1.498"
$ws.Range("E28").Value = "  +0.32%  "

Set-TextValue $ws.Range("D29") "4.118"
$ws.Range("E29").Value = "  +1.68%  "

Set-TextValue $ws.Range("D30") "4.075"
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("E31").Value = "  -0.27%  "

Set-TextValue $ws.Range("D32") "0.05373"
$ws.Range("E32").Value = "  +3.00%  "

Set-TextValue $ws.Range("D33") "1.859"
$ws.Range("E33").Value = "  +1.07%  "

Set-TextValue $ws.Range("D34") "0.7433"
$ws.Range("E34").Value = "  +1.54%  "

Set-TextValue $ws.Range("D35") "1.140"
$ws.Range("E35").Value = "  -1.23%  "

Set-TextValue $ws.Range("D36") "2.653"
$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("D37").Value = "1.244.81"
$ws.Range("E37").Value = "  -2.56%  "

Set-TextValue $ws.Range("D38") "0.01786"
$ws.Range("E38").Value = "  +0.27%  "

Set-TextValue $ws.Range("D39") "2.749"
$ws.Range("E39").Value = "  -0.06%  "

Set-TextValue $ws.Range("D40") "6.587"
$ws.Range("E40").Value = "  +3.37%  "

Set-TextValue $ws.Range("D41") "0.9046"
$ws.Range("E41").Value = "  +1.35%  "

Set-TextValue $ws.Range("D42") "0.9977"
$ws.Range("E42").Value = "  -0.27%  "

Set-TextValue $ws.Range("D43") "101.66"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "1.984.92"
$ws.Range("E44").Value = "  +0.30%  "

Set-TextValue $ws.Range("D45") "0.00000000122"
$ws.Range("E45").Value = "  +2.92%  "

Set-TextValue $ws.Range("D46") "64.82"
$ws.Range("E46").Value = "  +2.16%  "

Set-TextValue $ws.Range("D47") "0.5099"

Set-TextValue $ws.Range("D48") "0.4072"
$ws.Range("E48").Value = "  +2.71%  "

Set-TextValue $ws.Range("D49") "9.009"
$ws.Range("E49").Value = "  +1.52%  "

Set-TextValue $ws.Range("D50") "6.765"
$ws.Range("E50").Value = "  +1.44%  "

Set-TextValue $ws.Range("D51") "0.05763"
$ws.Range("E51").Value = "  +0.19%  "
